$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 338
$ws.Range("I42").Value = 338
$ws.Range("K42").Value = 1014
$ws.Range("M42").Value = -784
$ws.Range("H53").Value = 397.6
$ws.Range("I53").Value = 587.5
$ws.Range("J53").Value = 271
$ws.Range("K53").Value = 587.5
$ws.Range("L53").Value = 271
$ws.Range("M53").Value = 49.5
$ws.Range("N53").Value = -1545
$ws.Range("H80").Value = 30419.8
$ws.Range("I80").Value = 624.5
$ws.Range("K80").Value = 1873.5
$ws.Range("M80").Value = -875.5
$ws.Range("H83").Value = 30419.8
$ws.Range("I83").Value = 624.5
$ws.Range("K83").Value = 5620.5
$ws.Range("M83").Value = -628.5
$ws.Range("H132").Value = 2317.5
$ws.Range("I132").Value = 2313.6365
$ws.Range("K132").Value = 6940.9095
$ws.Range("M132").Value = -4410.9095

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2318.4
$ws.Range("I45").Value = 2318.4
$ws.Range("K45").Value = 2318.4
$ws.Range("M45").Value = -1941.4
$ws.Range("H61").Value = 1268.75
$ws.Range("I61").Value = 1268.75
$ws.Range("K61").Value = 1268.75
$ws.Range("M61").Value = -1056.75
$ws.Range("H74").Value = 2778.2727
$ws.Range("I74").Value = 866.8333
$ws.Range("J74").Value = 5072
$ws.Range("K74").Value = 866.8333
$ws.Range("L74").Value = 5072
$ws.Range("M74").Value = 7.166699999999992
$ws.Range("N74").Value = -6820
$ws.Range("H77").Value = 2778.2727
$ws.Range("I77").Value = 866.8333
$ws.Range("J77").Value = 5072
$ws.Range("K77").Value = 4334.1665
$ws.Range("L77").Value = 25360
$ws.Range("M77").Value = 33.83349999999973
$ws.Range("N77").Value = -34096
$ws.Range("H97").Value = 517.1667
$ws.Range("I97").Value = 525.625
$ws.Range("J97").Value = 449.5
$ws.Range("K97").Value = 525.625
$ws.Range("L97").Value = 449.5
$ws.Range("M97").Value = -29.625
$ws.Range("N97").Value = -1441.5
$ws.Range("H122").Value = 3798.8
$ws.Range("I122").Value = 2697.6
$ws.Range("J122").Value = 4900
$ws.Range("K122").Value = 8092.799999999999
$ws.Range("L122").Value = 14700
$ws.Range("M122").Value = -5642.799999999999
$ws.Range("N122").Value = -19600
$ws.Range("H132").Value = 2332.353
$ws.Range("I132").Value = 2179.4194
$ws.Range("K132").Value = 6538.2582
$ws.Range("M132").Value = -4008.2582
$ws.Range("H136").Value = 1268.75
$ws.Range("I136").Value = 1268.75
$ws.Range("K136").Value = 3806.25
$ws.Range("M136").Value = -1256.25
$ws.Range("H141").Value = 80000.5
$ws.Range("J141").Value = 80000.5
$ws.Range("L141").Value = 80000.5
$ws.Range("N141").Value = -90360.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1371.36
$ws.Range("I134").Value = 766.7143
$ws.Range("J134").Value = 4545.75
$ws.Range("K134").Value = 2300.1429
$ws.Range("L134").Value = 13637.25
$ws.Range("M134").Value = 234.8571000000002
$ws.Range("N134").Value = -18707.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2988.64
$ws.Range("I132").Value = 2517.1765
$ws.Range("K132").Value = 7551.529500000001
$ws.Range("M132").Value = -5021.529500000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1774.6
$ws.Range("I68").Value = 1485.6
$ws.Range("K68").Value = 4456.799999999999
$ws.Range("M68").Value = -3645.799999999999
$ws.Range("H71").Value = 1774.6
$ws.Range("I71").Value = 1485.6
$ws.Range("K71").Value = 13370.4
$ws.Range("M71").Value = -9314.4
$ws.Range("H75").Value = 442.4
$ws.Range("I75").Value = 442.4
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 1327.2
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -329.1999999999998
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 442.4
$ws.Range("I78").Value = 442.4
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 3981.6
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = 1010.4
$ws.Range("N78").ClearContents()
$ws.Range("H80").Value = 5858.6665
$ws.Range("I80").Value = 5998
$ws.Range("J80").Value = 5830.8
$ws.Range("K80").Value = 17994
$ws.Range("L80").Value = 17492.4
$ws.Range("M80").Value = -17058
$ws.Range("N80").Value = -19364.4
$ws.Range("H83").Value = 5858.6665
$ws.Range("I83").Value = 5998
$ws.Range("J83").Value = 5830.8
$ws.Range("K83").Value = 53982
$ws.Range("L83").Value = 52477.2
$ws.Range("M83").Value = -49302
$ws.Range("N83").Value = -61837.2
$ws.Range("H98").Value = 776
$ws.Range("J98").Value = 675
$ws.Range("L98").Value = 2025
$ws.Range("N98").Value = -5021
$ws.Range("H129").Value = 3742.7
$ws.Range("I129").Value = 3198.6667
$ws.Range("J129").Value = 3975.8572
$ws.Range("K129").Value = 9596.000100000001
$ws.Range("L129").Value = 11927.5716
$ws.Range("M129").Value = -4596.000100000001
$ws.Range("N129").Value = -21927.5716
$ws.Range("H132").Value = 2994.4546
$ws.Range("I132").Value = 963
$ws.Range("J132").Value = 6549.5
$ws.Range("K132").Value = 8667
$ws.Range("L132").Value = 58945.5
$ws.Range("M132").Value = -6137
$ws.Range("N132").Value = -64005.5
$ws.Range("H140").Value = 2359.9
$ws.Range("I140").Value = 2359.9
$ws.Range("K140").Value = 7079.700000000001
$ws.Range("M140").Value = -1899.700000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 14000
$ws.Range("J94").Value = 14000
$ws.Range("L94").Value = 14000
$ws.Range("N94").Value = -15352
$ws.Range("H102").Value = 1186.4286
$ws.Range("I102").Value = 526
$ws.Range("J102").Value = 3094.3333
$ws.Range("K102").Value = 526
$ws.Range("L102").Value = 3094.3333
$ws.Range("M102").Value = 1096
$ws.Range("N102").Value = -6338.3333
$ws.Range("H122").Value = 848786.9399999999
$ws.Range("I122").Value = 113137.89
$ws.Range("J122").Value = 2503997.2
$ws.Range("K122").Value = 339413.67
$ws.Range("L122").Value = 7511991.600000001
$ws.Range("M122").Value = -336963.67
$ws.Range("N122").Value = -7516891.600000001
$ws.Range("H126").Value = 3168.8462
$ws.Range("I126").Value = 1959
$ws.Range("J126").Value = 3925
$ws.Range("K126").Value = 5877
$ws.Range("L126").Value = 11775
$ws.Range("M126").Value = -3407
$ws.Range("N126").Value = -16715

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1646.3182
$ws.Range("I7").Value = 1772.4667
$ws.Range("J7").Value = 1376
$ws.Range("K7").Value = 1772.4667
$ws.Range("L7").Value = 1376
$ws.Range("M7").Value = -1660.4667
$ws.Range("N7").Value = -1600
$ws.Range("H16").Value = 6293.5713
$ws.Range("I16").Value = 5117.0625
$ws.Range("J16").Value = 10058.4
$ws.Range("K16").Value = 5117.0625
$ws.Range("L16").Value = 10058.4
$ws.Range("M16").Value = -4947.0625
$ws.Range("N16").Value = -10398.4
$ws.Range("H22").Value = 647.6667
$ws.Range("I22").Value = 754.8333
$ws.Range("K22").Value = 754.8333
$ws.Range("M22").Value = -459.8333
$ws.Range("H27").Value = 647.6667
$ws.Range("I27").Value = 754.8333
$ws.Range("K27").Value = 754.8333
$ws.Range("M27").Value = -647.8333
$ws.Range("H40").Value = 1818.1428
$ws.Range("I40").Value = 1818.1428
$ws.Range("K40").Value = 1818.1428
$ws.Range("M40").Value = -1682.1428
$ws.Range("H122").Value = 3425
$ws.Range("I122").Value = 3425
$ws.Range("K122").Value = 10275
$ws.Range("M122").Value = -7825
$ws.Range("H126").Value = 1646.3182
$ws.Range("I126").Value = 1772.4667
$ws.Range("J126").Value = 1376
$ws.Range("K126").Value = 5317.4001
$ws.Range("L126").Value = 4128
$ws.Range("M126").Value = -2847.4001
$ws.Range("N126").Value = -9068
$ws.Range("H136").Value = 3246.1667
$ws.Range("I136").Value = 2220.5
$ws.Range("K136").Value = 6661.5
$ws.Range("M136").Value = -4111.5
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 20998.25
$ws.Range("J15").Value = 20998.25
$ws.Range("L15").Value = 20998.25
$ws.Range("N15").Value = -21574.25
$ws.Range("H113").Value = 2289.762
$ws.Range("I113").Value = 1904.8889
$ws.Range("K113").Value = 5714.6667
$ws.Range("M113").Value = -3544.6667
$ws.Range("H132").Value = 1659.439
$ws.Range("I132").Value = 1245.8125
$ws.Range("K132").Value = 3737.4375
$ws.Range("M132").Value = -1207.4375
$ws.Range("H136").Value = 25671
$ws.Range("I136").Value = 1109.7037
$ws.Range("J136").Value = 73039.21000000001
$ws.Range("K136").Value = 3329.1111
$ws.Range("L136").Value = 219117.63
$ws.Range("M136").Value = -779.1111000000001
$ws.Range("N136").Value = -224217.63
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360
